$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 461-462; existing rows 461-557 shift down to 463-559.
$ws.Rows("461:462").Insert()

# --- New row 461 ---
$ws.Range("A461").Value2 = 6
$ws.Range("B461").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C461").Value2 = "Metropolitana"
$ws.Range("D461").Value2 = 44522
$ws.Range("E461").Value2 = 13
$ws.Range("F461").Value2 = 100112031
$ws.Range("G461").Value2 = "Poroto verde"
$ws.Range("H461").Value2 = "Magnum"
$ws.Range("I461").Value2 = "Primera"
$ws.Range("J461").Value2 = 400
$ws.Range("K461").Value2 = 33000
$ws.Range("L461").Value2 = 35000
$ws.Range("M461").Value2 = 34150
$ws.Range("N461").Value2 = "`$/saco 25 kilos"
$ws.Range("O461").Value2 = "Región de O'Higgins"
$ws.Range("P461").Value2 = 1366
$ws.Range("Q461").Value2 = 25
$ws.Range("R461").Value2 = "Hortaliza"

# --- New row 462 ---
$ws.Range("A462").Value2 = 6
$ws.Range("B462").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C462").Value2 = "Metropolitana"
$ws.Range("D462").Value2 = 44522
$ws.Range("E462").Value2 = 13
$ws.Range("F462").Value2 = 100112031
$ws.Range("G462").Value2 = "Poroto verde"
$ws.Range("H462").Value2 = "Sin especificar"
$ws.Range("I462").Value2 = "Primera"
$ws.Range("J462").Value2 = 250
$ws.Range("K462").Value2 = 50000
$ws.Range("L462").Value2 = 55000
$ws.Range("M462").Value2 = 52000
$ws.Range("N462").Value2 = "`$/malla 25 kilos"
$ws.Range("O462").Value2 = "Provincia del Elquí"
$ws.Range("P462").Value2 = 2080
$ws.Range("Q462").Value2 = 25
$ws.Range("R462").Value2 = "Hortaliza"
